$d = $word.ActiveDocument

# Merge the split runs of the Title paragraph into a single run.
$d.Content.Find.Execute("Answers: Introduction to Matrices", $false, $false, $false, $false, $false, $true, 1, $false, "Answers: Introduction to Matrices", 2) | Out-Null

# Merge the split runs of the Author paragraph into a single run.
$d.Content.Find.Execute("Jessica Taberner", $false, $false, $false, $false, $false, $true, 1, $false, "Jessica Taberner", 2) | Out-Null

# Merge the split runs of the Abstract paragraph into a single run.
$d.Content.Find.Execute("Answers to a selection of questions on matrices.", $false, $false, $false, $false, $false, $true, 1, $false, "Answers to a selection of questions on matrices.", 2) | Out-Null
